# Update crypto price/volume data as scraped on Fri Apr 14 11:02:39 UTC 2023
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D that would otherwise be auto-converted to numbers by Excel
# (losing their exact textual formatting) are forced to Text format first.
$textFormatCells = @("D5","D6","D7","D8","D9","D10","D12","D14","D15","D16","D19","D20","D21","D22","D24","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D48","D51")
foreach ($addr in $textFormatCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "30.816.71"
$ws.Range("E2").Value = "  +1.85%  "
$ws.Range("D3").Value = "2.116.89"
$ws.Range("E3").Value = "  +6.32%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "332.89"
$ws.Range("E5").Value = "  +2.51%  "
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").Value = "0.5350"
$ws.Range("E7").Value = "  +4.92%  "
$ws.Range("D8").Value = "0.4405"
$ws.Range("D9").Value = "0.09017"
$ws.Range("E9").Value = "  +3.30%  "
$ws.Range("D10").Value = "47.37"
$ws.Range("E10").Value = "  +11.03%  "
$ws.Range("E11").Value = "  +4.50%  "
$ws.Range("D12").Value = "24.99"
$ws.Range("E12").Value = "  +3.26%  "
$ws.Range("D13").Value = "2.114.07"
$ws.Range("E13").Value = "  +6.21%  "
$ws.Range("D14").Value = "6.771"
$ws.Range("E14").Value = "  +4.43%  "
$ws.Range("D15").Value = "7.821"
$ws.Range("E15").Value = "  +5.83%  "
$ws.Range("D16").Value = "96.99"
$ws.Range("E16").Value = "  +3.20%  "
$ws.Range("E17").Value = "  +0.15%  "
$ws.Range("E18").Value = "  +1.29%  "
$ws.Range("D19").Value = "0.06675"
$ws.Range("E19").Value = "  +1.90%  "
$ws.Range("D20").Value = "19.17"
$ws.Range("E20").Value = "  +1.73%  "
$ws.Range("D21").Value = "1.001"
$ws.Range("D22").Value = "6.349"
$ws.Range("E22").Value = "  +4.45%  "
$ws.Range("D23").Value = "30.882.16"
$ws.Range("E23").Value = "  +1.88%  "
$ws.Range("D24").Value = "12.34"
$ws.Range("E24").Value = "  +7.20%  "
$ws.Range("D25").Value = "2.363.47"
$ws.Range("E25").Value = "  +6.56%  "
$ws.Range("D26").Value = "2.294"
$ws.Range("E26").Value = "  +3.92%  "
$ws.Range("D27").Value = "22.79"
$ws.Range("E27").Value = "  +1.01%  "
$ws.Range("D28").Value = "2.601"
$ws.Range("E28").Value = "  +9.36%  "
$ws.Range("D29").Value = "163.38"
$ws.Range("E29").Value = "  +0.28%  "
$ws.Range("D30").Value = "133.43"
$ws.Range("E30").Value = "  +1.96%  "
$ws.Range("D31").Value = "1.180"
$ws.Range("E31").Value = "  +4.01%  "
$ws.Range("D32").Value = "0.1084"
$ws.Range("E32").Value = "  +2.79%  "
$ws.Range("D33").Value = "6.250"
$ws.Range("E33").Value = "  +3.31%  "
$ws.Range("D34").Value = "4.018"
$ws.Range("E34").Value = "  +5.51%  "
$ws.Range("D35").Value = "1.547"
$ws.Range("E35").Value = "  +17.03%  "
$ws.Range("D36").Value = "0.02604"
$ws.Range("E36").Value = "  +4.79%  "
$ws.Range("D37").Value = "5.564"
$ws.Range("E37").Value = "  +3.23%  "
$ws.Range("D38").Value = "12.90"
$ws.Range("E38").Value = "  +9.69%  "
$ws.Range("D39").Value = "0.06771"
$ws.Range("E39").Value = "  +4.20%  "
$ws.Range("D40").Value = "9.490"
$ws.Range("E40").Value = "  +6.69%  "
$ws.Range("D41").Value = "0.2290"
$ws.Range("E41").Value = "  +4.79%  "
$ws.Range("D42").Value = "0.6863"
$ws.Range("E42").Value = "  +4.29%  "
$ws.Range("D43").Value = "1.247"
$ws.Range("E43").Value = "  +2.05%  "
$ws.Range("D44").Value = "0.6461"
$ws.Range("E44").Value = "  +5.32%  "
$ws.Range("D45").Value = "14.16"
$ws.Range("E45").Value = "  +4.19%  "
$ws.Range("D46").Value = "1.000"
$ws.Range("E46").Value = "  -0.13%  "
$ws.Range("D47").Value = "2.230"
$ws.Range("E47").Value = "  +1.82%  "
$ws.Range("D48").Value = "3.657"
$ws.Range("E48").Value = "  -0.14%  "
$ws.Range("E49").Value = "  +4.39%  "
$ws.Range("E50").Value = "  +4.73%  "
$ws.Range("D51").Value = "121.87"
$ws.Range("E51").Value = "  -1.91%  "
